# 自动更新Excel文件 - 2025-12-07 23:11:57
# For each inventory row, decrement the "剩余" (remaining days) counter by one.
# When the remaining-day counter has run out (i.e. it was 1, meaning the
# supply was due today), the row is restocked: "剩余" is reset back to the
# full "总天" duration and "开始时间" (start date) is rolled forward to the
# new restock date (2025-12-08).
#
# Rows whose start-date value is not a well-formed 8-digit date (e.g. stray
# data-entry typos) are left untouched, since they could not be reliably
# re-evaluated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStartDate = 20251208

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $total = $ws.Cells.Item($r, 4).Value2      # D: 总天 (total days)
    $remaining = $ws.Cells.Item($r, 5).Value2  # E: 剩余 (remaining days)
    $startDate = $ws.Cells.Item($r, 6).Value2  # F: 开始时间 (start date)

    if ($total -eq $null -or $remaining -eq $null -or $startDate -eq $null) {
        continue
    }

    $startDateText = [string]([int]$startDate)
    if ($startDateText.Length -ne 8) {
        # Malformed / unexpected date value - skip this row untouched.
        continue
    }

    if ($remaining -eq 1) {
        $ws.Cells.Item($r, 5).Value = $total
        $ws.Cells.Item($r, 6).Value = $newStartDate
    } else {
        $ws.Cells.Item($r, 5).Value = $remaining - 1
    }
}
